$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with two new columns P and Q ---
# Copy the formatting of O1 (bold/border style) onto the new header cells.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Update data rows 2-25 ---
# Columns I, K, M, O swap their 1/2 pattern, and new columns P, Q (=2) are appended.
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}

Write-Output "done"
